$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to stage text values so Excel does not auto-coerce
# numeric-looking strings (e.g. "312.53") into real numbers when the
# target cells must stay plain text, matching the source data feed.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

function Set-TextValue {
    param($Range, [string]$Text)
    $helper.Value = $Text
    $helper.Copy()
    $Range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("D2") '27.989.73'
Set-TextValue $ws.Range("E2") '  -0.36%  '
Set-TextValue $ws.Range("D3") '1.859.26'
Set-TextValue $ws.Range("E3") '  -0.92%  '
Set-TextValue $ws.Range("E4") '  +0.12%  '
Set-TextValue $ws.Range("D5") '312.53'
Set-TextValue $ws.Range("E5") '  -0.32%  '
Set-TextValue $ws.Range("E6") '  +0.10%  '
Set-TextValue $ws.Range("D7") '0.5144'
Set-TextValue $ws.Range("E7") '  +1.41%  '
Set-TextValue $ws.Range("E8") '  -0.37%  '
Set-TextValue $ws.Range("D9") '0.08240'
Set-TextValue $ws.Range("E9") '  -8.24%  '
Set-TextValue $ws.Range("E10") '  -1.18%  '
Set-TextValue $ws.Range("D11") '41.47'
Set-TextValue $ws.Range("E11") '  -0.32%  '
Set-TextValue $ws.Range("D12") '6.187'
Set-TextValue $ws.Range("E12") '  -2.41%  '
Set-TextValue $ws.Range("D13") '20.56'
Set-TextValue $ws.Range("E13") '  -0.83%  '
Set-TextValue $ws.Range("D14") '1.866.31'
Set-TextValue $ws.Range("E14") '  -0.48%  '
Set-TextValue $ws.Range("D15") '7.251'
Set-TextValue $ws.Range("E15") '  +0.62%  '
Set-TextValue $ws.Range("E16") '  +0.07%  '
Set-TextValue $ws.Range("D17") '0.00001096'
Set-TextValue $ws.Range("E17") '  -0.99%  '
Set-TextValue $ws.Range("D18") '90.56'
Set-TextValue $ws.Range("E18") '  -0.68%  '
Set-TextValue $ws.Range("D19") '0.06652'
Set-TextValue $ws.Range("E19") '  +0.83%  '
Set-TextValue $ws.Range("D20") '17.67'
Set-TextValue $ws.Range("E20") '  -2.60%  '
Set-TextValue $ws.Range("D22") '6.002'
Set-TextValue $ws.Range("E22") '  -1.77%  '
Set-TextValue $ws.Range("D23") '28.019.33'
Set-TextValue $ws.Range("E23") '  -0.32%  '
Set-TextValue $ws.Range("D24") '11.06'
Set-TextValue $ws.Range("E24") '  -3.11%  '
Set-TextValue $ws.Range("D25") '2.247'
Set-TextValue $ws.Range("E25") '  -1.25%  '
Set-TextValue $ws.Range("D26") '2.075.35'
Set-TextValue $ws.Range("E26") '  -0.79%  '
Set-TextValue $ws.Range("D27") '2.515'
Set-TextValue $ws.Range("E27") '  -0.89%  '
Set-TextValue $ws.Range("D28") '158.10'
Set-TextValue $ws.Range("E28") '  +0.71%  '
Set-TextValue $ws.Range("D29") '20.44'
Set-TextValue $ws.Range("E29") '  -1.48%  '
Set-TextValue $ws.Range("D30") '124.60'
Set-TextValue $ws.Range("E30") '  -1.66%  '
Set-TextValue $ws.Range("E31") '  +1.20%  '
Set-TextValue $ws.Range("E32") '  -3.12%  '
Set-TextValue $ws.Range("D33") '5.987'
Set-TextValue $ws.Range("E33") '  +6.74%  '
Set-TextValue $ws.Range("D34") '3.601'
Set-TextValue $ws.Range("E34") '  -0.10%  '
Set-TextValue $ws.Range("D35") '9.321'
Set-TextValue $ws.Range("E35") '  -3.16%  '
Set-TextValue $ws.Range("D36") '0.02414'
Set-TextValue $ws.Range("E36") '  -0.18%  '
Set-TextValue $ws.Range("D37") '0.06490'
Set-TextValue $ws.Range("E37") '  -1.30%  '
Set-TextValue $ws.Range("D38") '0.2168'
Set-TextValue $ws.Range("D39") '0.6538'
Set-TextValue $ws.Range("E39") '  +2.22%  '
Set-TextValue $ws.Range("D40") '1.196'
Set-TextValue $ws.Range("E40") '  -0.79%  '
Set-TextValue $ws.Range("D41") '5.013'
Set-TextValue $ws.Range("E41") '  +1.96%  '
Set-TextValue $ws.Range("D42") '1.223'
Set-TextValue $ws.Range("E42") '  -3.62%  '
Set-TextValue $ws.Range("D43") '11.14'
Set-TextValue $ws.Range("E43") '  -2.93%  '
Set-TextValue $ws.Range("D44") '0.6135'
Set-TextValue $ws.Range("E44") '  +1.79%  '
Set-TextValue $ws.Range("D45") '12.96'
Set-TextValue $ws.Range("E45") '  -1.46%  '
Set-TextValue $ws.Range("D46") '1.280'
Set-TextValue $ws.Range("E46") '  +0.35%  '
Set-TextValue $ws.Range("D47") '3.664'
Set-TextValue $ws.Range("E47") '  -0.26%  '
Set-TextValue $ws.Range("E48") '  +0.62%  '
Set-TextValue $ws.Range("D49") '1.214'
Set-TextValue $ws.Range("E49") '  -1.86%  '
Set-TextValue $ws.Range("D50") '120.44'
Set-TextValue $ws.Range("E50") '  -0.69%  '
Set-TextValue $ws.Range("D51") '78.36'
Set-TextValue $ws.Range("E51") '  -1.72%  '

$helper.Clear()
